# Re-process the metadata sheet with the newly curated dimensions.
# - Several columns that were previously modeled as "dimension" concepts
#   (continente, area-nacionalidad, nivel-estudios-agregado, grado-de-formacion)
#   are now curated as "measure" concepts instead, which also changes their
#   "dim"/"medida" role (row 3) and their value-type ("skos:Concept" -> "xsd:int", row 4).
# - The "aragon" column is dropped as its own dimension and the "F" column now
#   reuses the refArea / URI-Comunidad mapping instead.
# - The obsolete mapping-file row (row 5) is removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: iaest-*:<concept> labels
$ws.Range("A2").Value = "iaest-measure:continente"
$ws.Range("C2").Value = "iaest-measure:area-nacionalidad"
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("G2").Value = "iaest-measure:nivel-estudios-agregado"

# Row 3: dim / medida role
$ws.Range("A3").Value = "medida"
$ws.Range("C3").Value = "medida"
$ws.Range("G3").Value = "medida"

# Row 4: value type / URI mapping
$ws.Range("A4").Value = "xsd:int"
$ws.Range("C4").Value = "xsd:int"
$ws.Range("F4").Value = "URI-Comunidad"
$ws.Range("G4").Value = "xsd:int"

# Row 5 (the mapping-*.xlsx file references) is no longer needed.
$ws.Rows.Item(5).Delete()
